$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 63
$prevRow = $row - 1

$ws.Cells.Item($row, 1).Value = 45967
$ws.Cells.Item($row, 2).Value = "22,1388"
$ws.Cells.Item($row, 3).Value = "16,1814"
$ws.Cells.Item($row, 4).Value = "15,4288"
$ws.Cells.Item($row, 5).Value = "15,4288"

# Match the date-column formatting used by the rows above (style of A62).
$ws.Cells.Item($prevRow, 1).Copy()
$ws.Cells.Item($row, 1).PasteSpecial(-4122)
$ws.Cells.Item($row, 1).Value = 45967
